$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be interpreted as text while we write the
# values (several prices look like plain numbers, e.g. "0.3740", and
# would otherwise be auto-converted to a Number and lose trailing zeros).
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '22.369.19'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '1.566.86'
$ws.Range('E3').Value = '  +0.44%  '
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').Value = '1.005'
$ws.Range('E5').Value = '  +0.37%  '
$ws.Range('D6').Value = '289.06'
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('D7').Value = '0.3740'
$ws.Range('E7').Value = '  +0.60%  '
$ws.Range('D8').Value = '49.19'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = '0.3373'
$ws.Range('E9').Value = '  -0.98%  '
$ws.Range('D10').Value = '0.07438'
$ws.Range('E10').Value = '  -2.45%  '
$ws.Range('E11').Value = '  -3.86%  '
$ws.Range('E12').Value = '  +0.64%  '
$ws.Range('D13').Value = '20.73'
$ws.Range('E13').Value = '  -2.95%  '
$ws.Range('D14').Value = '5.869'
$ws.Range('E14').Value = '  -2.60%  '
$ws.Range('D15').Value = '6.849'
$ws.Range('E15').Value = '  -0.89%  '
$ws.Range('D16').Value = '1.563.75'
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('E17').Value = '  -1.63%  '
$ws.Range('D18').Value = '88.90'
$ws.Range('E18').Value = '  -0.95%  '
$ws.Range('D19').Value = '0.06698'
$ws.Range('E19').Value = '  -0.34%  '
$ws.Range('E20').Value = '  +0.51%  '
$ws.Range('D21').Value = '6.128'
$ws.Range('E21').Value = '  -1.47%  '
$ws.Range('D22').Value = '16.20'
$ws.Range('E22').Value = '  -1.95%  '
$ws.Range('D23').Value = '11.81'
$ws.Range('E23').Value = '  -1.13%  '
$ws.Range('D24').Value = '22.355.10'
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('D25').Value = '2.370'
$ws.Range('E25').Value = '  -1.41%  '
$ws.Range('D26').Value = '2.511'
$ws.Range('E26').Value = '  -10.80%  '
$ws.Range('E27').Value = '  -1.36%  '
$ws.Range('D28').Value = '147.14'
$ws.Range('E28').Value = '  +0.75%  '
$ws.Range('E29').Value = '  +0.36%  '
$ws.Range('D30').Value = '124.88'
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('D31').Value = '1.737.45'
$ws.Range('E31').Value = '  +0.21%  '
$ws.Range('D32').Value = '1.993'
$ws.Range('E32').Value = '  -0.64%  '
$ws.Range('D33').Value = '0.9799'
$ws.Range('E33').Value = '  -2.51%  '
$ws.Range('D34').Value = '5.891'
$ws.Range('E34').Value = '  -4.40%  '
$ws.Range('D35').Value = '9.690'
$ws.Range('E35').Value = '  -3.17%  '
$ws.Range('D36').Value = '0.08412'
$ws.Range('E36').Value = '  -0.81%  '
$ws.Range('D37').Value = '1.374'
$ws.Range('E37').Value = '  +4.20%  '
$ws.Range('E38').Value = '  -3.54%  '
$ws.Range('E39').Value = '  -2.31%  '
$ws.Range('E40').Value = '  -0.50%  '
$ws.Range('D41').Value = '5.339'
$ws.Range('E41').Value = '  -2.96%  '
$ws.Range('E42').Value = '  -2.64%  '
$ws.Range('D43').Value = '10.91'
$ws.Range('E43').Value = '  -6.38%  '
$ws.Range('E44').Value = '  +0.40%  '
$ws.Range('D45').Value = '13.74'
$ws.Range('E45').Value = '  -2.51%  '
$ws.Range('D46').Value = '3.776'
$ws.Range('E46').Value = '  +0.55%  '
$ws.Range('D47').Value = '0.5748'
$ws.Range('E47').Value = '  -3.53%  '
$ws.Range('D48').Value = '2.032'
$ws.Range('E48').Value = '  -2.59%  '
$ws.Range('D49').Value = '125.01'
$ws.Range('E49').Value = '  +0.70%  '
$ws.Range('D50').Value = '1.227'
$ws.Range('E50').Value = '  -2.94%  '
$ws.Range('D51').Value = '0.07299'
$ws.Range('E51').Value = '  +0.52%  '

# Restore the original (default) cell style on the Price column so only
# the text content changes, matching the source data.
$ws.Range('D2:D51').Style = 'Normal'
